# Comercializadora del Agro de Limarí - Palta
# Adds a new week of reported prices (week of 2023-12-07, serial 45267) as
# six new rows inserted before the existing row 861, pushing all the
# subsequent historical rows down by six rows (861->867 ... 899->905).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 861:899(+) down by 6 rows, inserting 6 blank rows at 861.
$ws.Range("A861:T866").Insert()

# Fixed values shared by every data row in this sheet.
$mercadoId = 2
$mercado   = "Comercializadora del Agro de Limarí"
$region    = "Coquimbo"
$codreg    = 4
$tipo      = "Fruta"
$productoId = 100106
$producto  = "Oleaginosos"
$categoriaId = 100106002
$categoria = "Palta"
$unidad    = "`$/kilo (en caja de 17 kilos)"
$origen    = "Provincia de Limarí"
$kgUnidad  = 1

# New rows to insert at 861..866 (D, K, L, M, N, O, P, S).
$newRows = @(
    @(45267, "Edranol", "Especial", 200, 2600, 2700, 2650, 2650),
    @(45267, "Edranol", "Primera",  240, 2300, 2400, 2350, 2350),
    @(45267, "Edranol", "Segunda",  200, 1900, 2000, 1950, 1950),
    @(45267, "Hass",    "Especial", 200, 3100, 3200, 3150, 3150),
    @(45267, "Hass",    "Primera",  300, 2800, 2900, 2850, 2850),
    @(45267, "Hass",    "Segunda",  240, 2500, 2600, 2550, 2550)
)

for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = 861 + $i
    $data = $newRows[$i]

    $ws.Cells.Item($r, 1).Value = $mercadoId
    $ws.Cells.Item($r, 2).Value = $mercado
    $ws.Cells.Item($r, 3).Value = $region
    $ws.Cells.Item($r, 4).Value = $data[0]
    $ws.Cells.Item($r, 5).Value = $codreg
    $ws.Cells.Item($r, 6).Value = $tipo
    $ws.Cells.Item($r, 7).Value = $productoId
    $ws.Cells.Item($r, 8).Value = $producto
    $ws.Cells.Item($r, 9).Value = $categoriaId
    $ws.Cells.Item($r, 10).Value = $categoria
    $ws.Cells.Item($r, 11).Value = $data[1]
    $ws.Cells.Item($r, 12).Value = $data[2]
    $ws.Cells.Item($r, 13).Value = $data[3]
    $ws.Cells.Item($r, 14).Value = $data[4]
    $ws.Cells.Item($r, 15).Value = $data[5]
    $ws.Cells.Item($r, 16).Value = $data[6]
    $ws.Cells.Item($r, 17).Value = $unidad
    $ws.Cells.Item($r, 18).Value = $origen
    $ws.Cells.Item($r, 19).Value = $data[7]
    $ws.Cells.Item($r, 20).Value = $kgUnidad
}
